$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new daily progress values (previously blank cells)
$ws.Range("P8").Value = 2
$ws.Range("N9").Value = 2
$ws.Range("N11").Value = 3
$ws.Range("N13").Value = 2

# Update the currently selected cell on the sheet
$ws.Range("J10").Select()

$wb.Save()
